# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# 1) Update the "last updated" date on the About sheet (C1) to 3/28/2024.
# 2) On RAF-capacity, raise the RAF for the two hydrogen technologies
#    (hydrogen combustion turbine, hydrogen combined cycle) from 0.3 to 1.

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("C1").Value = (Get-Date -Year 2024 -Month 3 -Day 28)

$capacitySheet = $wb.Worksheets.Item("RAF-capacity")
$capacitySheet.Range("B24").Value = 1
$capacitySheet.Range("B25").Value = 1

# Reflect the sheet that was active/selected when the workbook was saved.
$capacitySheet.Select()
$capacitySheet.Range("B25").Select()
